{"js": "// Update the 15 lattice-multiplication exercise cells (5 rows x 3 cols)\n// in the single table of the document with a new set of problems, while\n// preserving each cell's existing paragraph/run formatting (sz=32).\n//\n// Each cell's text is 5 lines (separated by manual line breaks):\n//   \"{A} x {B}\"\n//   \"  {B[0]}    {B[1]}\"\n//   \"  ----\"\n//   \"{A[0]}|    |\"\n//   \"{A[1]}|    |\"\n// where A and B are the two (two-digit) factors of the problem.\n\nconst newProblems = [\n  [16, 76], [65, 73], [21, 42],\n  [29, 96], [35, 64], [71, 84],\n  [90, 61], [38, 86], [34, 67],\n  [18, 29], [10, 39], [62, 57],\n  [81, 17], [53, 35], [14, 16],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\ntable.rows.items[0].cells.load(\"items\");\nawait context.sync();\nconst colCount = table.rows.items[0].cells.items.length;\n\n// First pass: fetch the first paragraph of every cell's body.\nconst paras = [];\nfor (let r = 0; r < rowCount; r++) {\n  const rowParas = [];\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const body = cell.body;\n    body.paragraphs.load(\"items\");\n    rowParas.push(body.paragraphs);\n  }\n  paras.push(rowParas);\n}\nawait context.sync();\n\n// Second pass: replace the text of each cell's first paragraph in place\n// (keeps the run's existing rPr, e.g. sz=32).\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const idx = r * colCount + c;\n    const [a, b] = newProblems[idx];\n    const aStr = String(a);\n    const bStr = String(b);\n    const line1 = `${a} x ${b}`;\n    const line2 = `  ${bStr[0]}    ${bStr[1]}`;\n    const line3 = \"  ----\";\n    const line4 = `${aStr[0]}|    |`;\n    const line5 = `${aStr[1]}|    |`;\n    const text = [line1, line2, line3, line4, line5].join(\"\\u000b\");\n\n    const para = paras[r][c].items[0];\n    para.insertText(text, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Update the 15 lattice-multiplication exercise cells (5 rows x 3 cols)\n# in the single table of the document with a new set of problems, while\n# preserving each cell's existing paragraph/run formatting (sz=32).\n#\n# Each cell's text is 5 lines (separated by manual line breaks, char 11 /\n# vertical-tab, which Word COM uses for \"<w:br/>\"):\n#   \"{A} x {B}\"\n#   \"  {B[0]}    {B[1]}\"\n#   \"  ----\"\n#   \"{A[0]}|    |\"\n#   \"{A[1]}|    |\"\n# where A and B are the two (two-digit) factors of the problem.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newProblems = @(\n    @(16, 76), @(65, 73), @(21, 42),\n    @(29, 96), @(35, 64), @(71, 84),\n    @(90, 61), @(38, 86), @(34, 67),\n    @(18, 29), @(10, 39), @(62, 57),\n    @(81, 17), @(53, 35), @(14, 16)\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n$br = [char]11\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $idx = (($r - 1) * $colCount) + ($c - 1)\n        $pair = $newProblems[$idx]\n        $a = $pair[0]\n        $b = $pair[1]\n        $aStr = [string]$a\n        $bStr = [string]$b\n        $aDigit0 = $aStr.Substring(0,1)\n        $aDigit1 = $aStr.Substring(1,1)\n        $bDigit0 = $bStr.Substring(0,1)\n        $bDigit1 = $bStr.Substring(1,1)\n\n        $line1 = \"$a x $b\"\n        $line2 = \"  $bDigit0    $bDigit1\"\n        $line3 = \"  ----\"\n        $line4 = \"$aDigit0|    |\"\n        $line5 = \"$aDigit1|    |\"\n\n        $text = \"$line1$br$line2$br$line3$br$line4$br$line5\"\n\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $text\n    }\n}\n"}
